$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new quarter columns N, O, P ---
$ws.Range("N1").Value = "31/12/2023"
$ws.Range("O1").Value = "31/03/2024"
$ws.Range("P1").Value = "30/06/2024"
# Copy header style (bold, centered, bordered) from M1 onto the new header cells
$ws.Range("M1").Copy()
$ws.Range("N1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows ---
$ws.Cells.Item(2, 14).Value = 5766859.776
$ws.Cells.Item(2, 15).Value = 5944329.216
$ws.Cells.Item(2, 16).Value = 6259262.976
$ws.Cells.Item(3, 14).Value = 945920
$ws.Cells.Item(3, 15).Value = 1025611.008
$ws.Cells.Item(3, 16).Value = 983409.024
$ws.Cells.Item(4, 14).Value = 233350
$ws.Cells.Item(4, 15).Value = 346931.008
$ws.Cells.Item(4, 16).Value = 281335.008
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(6, 14).Value = 318047.008
$ws.Cells.Item(6, 15).Value = 326055.008
$ws.Cells.Item(6, 16).Value = 359548.992
$ws.Cells.Item(7, 14).Value = 87974
$ws.Cells.Item(7, 15).Value = 87201
$ws.Cells.Item(7, 16).Value = 98012
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(9, 14).Value = 230968
$ws.Cells.Item(9, 15).Value = 185144
$ws.Cells.Item(9, 16).Value = 157303.008
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(11, 14).Value = 75581
$ws.Cells.Item(11, 15).Value = 80280
$ws.Cells.Item(11, 16).Value = 87210
$ws.Cells.Item(12, 14).Value = 288775.008
$ws.Cells.Item(12, 15).Value = 319337.984
$ws.Cells.Item(12, 16).Value = 345900.992
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(16, 14).Value = 7488
$ws.Cells.Item(16, 15).Value = 7683
$ws.Cells.Item(16, 16).Value = 7723
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(19, 14).Value = 110512
$ws.Cells.Item(19, 15).Value = 115957
$ws.Cells.Item(19, 16).Value = 128270
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(21, 14).Value = 55646
$ws.Cells.Item(21, 15).Value = 58290
$ws.Cells.Item(21, 16).Value = 64839
$ws.Cells.Item(22, 14).Value = 465172.992
$ws.Cells.Item(22, 15).Value = 483040.992
$ws.Cells.Item(22, 16).Value = 531748.992
$ws.Cells.Item(23, 14).Value = 3934074.112
$ws.Cells.Item(23, 15).Value = 3983280.896
$ws.Cells.Item(23, 16).Value = 4255911.936
$ws.Cells.Item(24, 14).Value = 132918
$ws.Cells.Item(24, 15).Value = 133058
$ws.Cells.Item(24, 16).Value = 142292
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 15).Value = 0
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(26, 14).Value = 5766859.776
$ws.Cells.Item(26, 15).Value = 5944329.216
$ws.Cells.Item(26, 16).Value = 6259262.976
$ws.Cells.Item(27, 14).Value = 876772.992
$ws.Cells.Item(27, 15).Value = 897742.976
$ws.Cells.Item(27, 16).Value = 1023318.976
$ws.Cells.Item(28, 14).Value = 122382
$ws.Cells.Item(28, 15).Value = 143467.008
$ws.Cells.Item(28, 16).Value = 110961
$ws.Cells.Item(29, 14).Value = 157164
$ws.Cells.Item(29, 15).Value = 147439.008
$ws.Cells.Item(29, 16).Value = 187646
$ws.Cells.Item(30, 14).Value = 52443
$ws.Cells.Item(30, 15).Value = 49623
$ws.Cells.Item(30, 16).Value = 50382
$ws.Cells.Item(31, 14).Value = 482380.992
$ws.Cells.Item(31, 15).Value = 498974.016
$ws.Cells.Item(31, 16).Value = 585577.984
$ws.Cells.Item(32, 14).Value = 0
$ws.Cells.Item(32, 15).Value = 0
$ws.Cells.Item(32, 16).Value = 0
$ws.Cells.Item(33, 14).Value = 0
$ws.Cells.Item(33, 15).Value = 0
$ws.Cells.Item(33, 16).Value = 0
$ws.Cells.Item(34, 14).Value = 62403
$ws.Cells.Item(34, 15).Value = 58240
$ws.Cells.Item(34, 16).Value = 88752
$ws.Cells.Item(35, 14).Value = 0
$ws.Cells.Item(35, 15).Value = 0
$ws.Cells.Item(35, 16).Value = 0
$ws.Cells.Item(36, 14).Value = 0
$ws.Cells.Item(36, 15).Value = 0
$ws.Cells.Item(36, 16).Value = 0
$ws.Cells.Item(37, 14).Value = 2536230.912
$ws.Cells.Item(37, 15).Value = 2590360.064
$ws.Cells.Item(37, 16).Value = 2689921.024
$ws.Cells.Item(38, 14).Value = 2173008.896
$ws.Cells.Item(38, 15).Value = 2205829.12
$ws.Cells.Item(38, 16).Value = 2209701.12
$ws.Cells.Item(39, 14).Value = 295
$ws.Cells.Item(39, 15).Value = 37
$ws.Cells.Item(39, 16).Value = 50
$ws.Cells.Item(40, 14).Value = 9908
$ws.Cells.Item(40, 15).Value = 10135
$ws.Cells.Item(40, 16).Value = 10362
$ws.Cells.Item(41, 14).Value = 317568.992
$ws.Cells.Item(41, 15).Value = 339308
$ws.Cells.Item(41, 16).Value = 429503.008
$ws.Cells.Item(42, 14).Value = 0
$ws.Cells.Item(42, 15).Value = 0
$ws.Cells.Item(42, 16).Value = 0
$ws.Cells.Item(43, 14).Value = 35450
$ws.Cells.Item(43, 15).Value = 35051
$ws.Cells.Item(43, 16).Value = 40305
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(44, 15).Value = 0
$ws.Cells.Item(44, 16).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(45, 16).Value = 0
$ws.Cells.Item(46, 14).Value = 444
$ws.Cells.Item(46, 15).Value = 1357
$ws.Cells.Item(46, 16).Value = 2445
$ws.Cells.Item(47, 14).Value = 2353412
$ws.Cells.Item(47, 15).Value = 2454869.048
$ws.Cells.Item(47, 16).Value = 2543577.912
$ws.Cells.Item(48, 14).Value = 351670.016
$ws.Cells.Item(48, 15).Value = 351670.016
$ws.Cells.Item(48, 16).Value = 351670.016
$ws.Cells.Item(49, 14).Value = 39499
$ws.Cells.Item(49, 15).Value = 39794
$ws.Cells.Item(49, 16).Value = 4895
$ws.Cells.Item(50, 14).Value = 0
$ws.Cells.Item(50, 15).Value = 0
$ws.Cells.Item(50, 16).Value = 0
$ws.Cells.Item(51, 14).Value = 1360557.952
$ws.Cells.Item(51, 15).Value = 1301362.048
$ws.Cells.Item(51, 16).Value = 1260851.968
$ws.Cells.Item(52, 14).Value = 0
$ws.Cells.Item(52, 15).Value = 104345
$ws.Cells.Item(52, 16).Value = 72116
$ws.Cells.Item(53, 14).Value = 0
$ws.Cells.Item(53, 15).Value = 0
$ws.Cells.Item(53, 16).Value = 0
$ws.Cells.Item(54, 14).Value = 0
$ws.Cells.Item(54, 15).Value = 0
$ws.Cells.Item(54, 16).Value = 0
$ws.Cells.Item(55, 14).Value = 601684.992
$ws.Cells.Item(55, 15).Value = 657697.984
$ws.Cells.Item(55, 16).Value = 854044.992
$ws.Cells.Item(56, 14).Value = 0
$ws.Cells.Item(56, 15).Value = 0
$ws.Cells.Item(56, 16).Value = 0
# Row 57: blank label row -> create empty N57:P57 matching style of M57
$ws.Range("M57").Copy()
$ws.Range("N57:P57").PasteSpecial(-4122)
# Row 58: blank label row -> create empty N58:P58 matching style of M58
$ws.Range("M58").Copy()
$ws.Range("N58:P58").PasteSpecial(-4122)
$ws.Cells.Item(59, 14).Value = 644823.104
$ws.Cells.Item(59, 15).Value = 640934.016
$ws.Cells.Item(59, 16).Value = 693891.008
$ws.Cells.Item(60, 14).Value = -347979.968
$ws.Cells.Item(60, 15).Value = -343764
$ws.Cells.Item(60, 16).Value = -375991.008
$ws.Cells.Item(61, 14).Value = 296842.976
$ws.Cells.Item(61, 15).Value = 297169.984
$ws.Cells.Item(61, 16).Value = 317900
$ws.Cells.Item(62, 14).Value = -10269
$ws.Cells.Item(62, 15).Value = -3441
$ws.Cells.Item(62, 16).Value = -4204
$ws.Cells.Item(63, 14).Value = -121690
$ws.Cells.Item(63, 15).Value = -114004
$ws.Cells.Item(63, 16).Value = -123682
$ws.Cells.Item(64, 14).Value = 0
$ws.Cells.Item(64, 15).Value = 0
$ws.Cells.Item(64, 16).Value = 0
$ws.Cells.Item(65, 14).Value = 5660
$ws.Cells.Item(65, 15).Value = 11194
$ws.Cells.Item(65, 16).Value = 6407
$ws.Cells.Item(66, 14).Value = 682
$ws.Cells.Item(66, 15).Value = -5012
$ws.Cells.Item(66, 16).Value = -7634
$ws.Cells.Item(67, 14).Value = 1165
$ws.Cells.Item(67, 15).Value = 3598
$ws.Cells.Item(67, 16).Value = -3315
$ws.Cells.Item(68, 14).Value = -26102
$ws.Cells.Item(68, 15).Value = -28809
$ws.Cells.Item(68, 16).Value = -56577
$ws.Cells.Item(69, 14).Value = 14703
$ws.Cells.Item(69, 15).Value = 16721
$ws.Cells.Item(69, 16).Value = 8783
$ws.Cells.Item(70, 14).Value = -40805
$ws.Cells.Item(70, 15).Value = -45530
$ws.Cells.Item(70, 16).Value = -65360
# Row 71: blank label row -> create empty N71:P71 matching style of M71
$ws.Range("M71").Copy()
$ws.Range("N71:P71").PasteSpecial(-4122)
# Row 72: blank label row -> create empty N72:P72 matching style of M72
$ws.Range("M72").Copy()
$ws.Range("N72:P72").PasteSpecial(-4122)
# Row 73: blank label row -> create empty N73:P73 matching style of M73
$ws.Range("M73").Copy()
$ws.Range("N73:P73").PasteSpecial(-4122)
$ws.Cells.Item(74, 14).Value = 146288.992
$ws.Cells.Item(74, 15).Value = 160696
$ws.Cells.Item(74, 16).Value = 128895
$ws.Cells.Item(75, 14).Value = 53811
$ws.Cells.Item(75, 15).Value = -42373
$ws.Cells.Item(75, 16).Value = -24049
$ws.Cells.Item(76, 14).Value = -86626
$ws.Cells.Item(76, 15).Value = -13065
$ws.Cells.Item(76, 16).Value = -59986
# Row 77: blank label row -> create empty N77:P77 matching style of M77
$ws.Range("M77").Copy()
$ws.Range("N77:P77").PasteSpecial(-4122)
# Row 78: blank label row -> create empty N78:P78 matching style of M78
$ws.Range("M78").Copy()
$ws.Range("N78:P78").PasteSpecial(-4122)
$ws.Cells.Item(79, 14).Value = -1002
$ws.Cells.Item(79, 15).Value = -913
$ws.Cells.Item(79, 16).Value = -1088
$ws.Cells.Item(80, 14).Value = 112472.016
$ws.Cells.Item(80, 15).Value = 104345
$ws.Cells.Item(80, 16).Value = 43772

$excel.CutCopyMode = 0
Write-Output "OK"